$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) storage for numeric-looking values in columns D and E so
# Excel does not silently convert them to actual numbers (losing formats like
# trailing zeros, thousand-dot separators, percent signs and padding spaces).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "57.939.56"
$ws.Range("E2").Value = "  +1.81%  "

$ws.Range("D3").Value = "2.347.76"
$ws.Range("E3").Value = "  +1.29%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "539.31"
$ws.Range("E5").Value = "  +1.78%  "

$ws.Range("D6").Value = "135.82"
$ws.Range("E6").Value = "  +2.59%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  +5.60%  "

$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").Value = "5.54"
$ws.Range("E10").Value = "  +4.63%  "

$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("E12").Value = "  +2.02%  "

$ws.Range("D13").Value = "23.77"
$ws.Range("E13").Value = "  +1.46%  "

$ws.Range("D14").Value = "2.765.85"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("D15").Value = "57.915.61"
$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "2.366.90"
$ws.Range("E17").Value = "  +1.19%  "

$ws.Range("D18").Value = "10.68"
$ws.Range("E18").Value = "  +2.61%  "

$ws.Range("D19").Value = "332.23"
$ws.Range("E19").Value = "  -1.14%  "

$ws.Range("E20").Value = "  +2.76%  "

$ws.Range("E21").Value = "  -0.95%  "

$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("E23").Value = "  +2.11%  "

$ws.Range("D24").Value = "0.166"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("E25").Value = "  -2.20%  "

$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "172.02"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "1.75"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +1.44%  "

$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  +0.51%  "

$ws.Range("E32").Value = "  +10.84%  "

$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.91%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.22"
$ws.Range("E36").Value = "  +5.94%  "

$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("D38").Value = "1.63"
$ws.Range("E38").Value = "  +4.17%  "

$ws.Range("D39").Value = "39.23"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("D40").Value = "145.47"
$ws.Range("E40").Value = "  -2.50%  "

$ws.Range("D41").Value = "293.21"
$ws.Range("E41").Value = "  +3.59%  "

$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("E43").Value = "  +1.18%  "

$ws.Range("D44").Value = "0.0949"
$ws.Range("E44").Value = "  +1.98%  "

$ws.Range("D45").Value = "19.26"
$ws.Range("E45").Value = "  +2.25%  "

$ws.Range("D46").Value = "0.0503"
$ws.Range("E46").Value = "  +0.45%  "

$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").Value = "0.386"
$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "17.46"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").Value = "11.07"
$ws.Range("E51").Value = "  +0.46%  "

